$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F header/stat labels and values (entered in the order that reproduces
# the shared-string table order seen in the target workbook)
$ws.Range("F2").Value = "Anzahl von Tests:"
$ws.Range("G2").Value = 10

$ws.Range("F3").Value = "Anzahl von fehlgeschlagene Tests"
$ws.Range("G3").Value = 6

$ws.Range("F4").Value = "Anzahl von gultige Tests"
$ws.Range("G4").Value = 4

$ws.Range("F1").Value = "Statistik"

$ws.Range("F6").Value = "Teststatistik nacher"
$ws.Range("G6").Value = 20

$ws.Range("F7").Value = "Anzahl von fehlgeschlagene Tests"
$ws.Range("G7").Value = 1

$ws.Range("F8").Value = "Anzahl von gute Tests "
$ws.Range("G8").Value = 19

# Column F width (matches the width recorded after authoring the new column)
$ws.Columns.Item(6).ColumnWidth = 30.8

# Restore the selection cell noted in the saved workbook
$ws.Range("F9").Select() | Out-Null
